$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "R.F.C.: THDRHRGHSRH" "R.F.C.: GODE561231GR8"
Replace-Text "CURP: RETRRERERTTRERER" "CURP: GODE560912HDFLNS09"
Replace-Text "Calle: COLONIA" "Calle: Unidad Volkswagen 1"
Replace-Text "Número Exterior: COLONIA" "Número Exterior: Sin colonia 2"
Replace-Text "Estado: PUEBLA" "Estado: Puebla"
Replace-Text "Ciudad: PUEBLAYORK" "Ciudad: Heroica Puebla de Zaragoza"
Replace-Text "País: MEXICO" "País: México"
Replace-Text "Puesto: DESARROLLADOR" "Puesto: Desarrollador"
Replace-Text "Correo Electrónico: juanito@gmail.com" "Correo Electrónico: juanpro@gmail.com"
Replace-Text "el puesto de DESARROLLADOR en Loma Expertos" "el puesto de Desarrollador en Loma Expertos"
